# feat: add 3 pc map
#
# Converts a block of cells (previously shared-string references to the
# "farming" string, text-typed but holding the digit "1") into plain
# numeric cells with value 1, matching the rest of the "pc map" grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each of these ranges currently holds text cells (t="s") pointing at the
# shared string "farming" (index 1 - which happens to render as "1").
# Re-assigning a numeric 1 turns them into normal numeric cells, same as
# their neighbours.
$ranges = @(
    "CE2:CM2",
    "CD3:CM3",
    "BZ4:CI4",
    "BZ5:CI5",
    "BZ6:CF6",
    "BZ7:CD7"
)

foreach ($rangeAddress in $ranges) {
    $ws.Range($rangeAddress).Value = 1
}

# Update the visible window / selection to match the edited area
# (scroll the viewport so AM1 is the top-left visible cell, then select
# the new active cell).
$win = $excel.ActiveWindow
$win.ScrollColumn = 39
$win.ScrollRow = 1
$ws.Range("CD15").Select()
